$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44294
$ws.Cells.Item(2, 8).Value = 'Magnum'
$ws.Cells.Item(2, 10).Value = 100
$ws.Cells.Item(2, 11).Value = 24000
$ws.Cells.Item(2, 12).Value = 25000
$ws.Cells.Item(2, 13).Value = 24500
$ws.Cells.Item(2, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(2, 15).Value = 'Región del Maule'
$ws.Cells.Item(2, 16).Value = 980

# Row 3
$ws.Cells.Item(3, 4).Value = 44160
$ws.Cells.Item(3, 8).Value = 'Magnum'
$ws.Cells.Item(3, 10).Value = 100
$ws.Cells.Item(3, 11).Value = 28000
$ws.Cells.Item(3, 12).Value = 30000
$ws.Cells.Item(3, 13).Value = 29000
$ws.Cells.Item(3, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(3, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(3, 16).Value = 1160

# Row 4
$ws.Cells.Item(4, 4).Value = 44461
$ws.Cells.Item(4, 8).Value = 'Sin especificar'
$ws.Cells.Item(4, 10).Value = 100
$ws.Cells.Item(4, 11).Value = 33000
$ws.Cells.Item(4, 12).Value = 34000
$ws.Cells.Item(4, 13).Value = 33500
$ws.Cells.Item(4, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(4, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(4, 16).Value = 1340

# Row 5
$ws.Cells.Item(5, 4).Value = 44253
$ws.Cells.Item(5, 8).Value = 'Magnum'
$ws.Cells.Item(5, 10).Value = 200
$ws.Cells.Item(5, 11).Value = 25000
$ws.Cells.Item(5, 12).Value = 26000
$ws.Cells.Item(5, 13).Value = 25500
$ws.Cells.Item(5, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(5, 15).Value = 'Región del Maule'
$ws.Cells.Item(5, 16).Value = 1020

# Row 6
$ws.Cells.Item(6, 4).Value = 44188
$ws.Cells.Item(6, 8).Value = 'Magnum'
$ws.Cells.Item(6, 10).Value = 100
$ws.Cells.Item(6, 11).Value = 38000
$ws.Cells.Item(6, 12).Value = 40000
$ws.Cells.Item(6, 13).Value = 39000
$ws.Cells.Item(6, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(6, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(6, 16).Value = 1560

# Row 7
$ws.Cells.Item(7, 4).Value = 44230
$ws.Cells.Item(7, 8).Value = 'Magnum'
$ws.Cells.Item(7, 10).Value = 100
$ws.Cells.Item(7, 11).Value = 22000
$ws.Cells.Item(7, 12).Value = 24000
$ws.Cells.Item(7, 13).Value = 23000
$ws.Cells.Item(7, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(7, 15).Value = 'Región del Maule'
$ws.Cells.Item(7, 16).Value = 920

# Row 8
$ws.Cells.Item(8, 4).Value = 44433
$ws.Cells.Item(8, 8).Value = 'Magnum'
$ws.Cells.Item(8, 10).Value = 100
$ws.Cells.Item(8, 11).Value = 25000
$ws.Cells.Item(8, 12).Value = 26000
$ws.Cells.Item(8, 13).Value = 25500
$ws.Cells.Item(8, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(8, 15).Value = 'Perú'
$ws.Cells.Item(8, 16).Value = 1020

# Row 9
$ws.Cells.Item(9, 4).Value = 44335
$ws.Cells.Item(9, 8).Value = 'Magnum'
$ws.Cells.Item(9, 10).Value = 100
$ws.Cells.Item(9, 11).Value = 35000
$ws.Cells.Item(9, 12).Value = 36000
$ws.Cells.Item(9, 13).Value = 35500
$ws.Cells.Item(9, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(9, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(9, 16).Value = 1420

# Row 10
$ws.Cells.Item(10, 4).Value = 44489
$ws.Cells.Item(10, 8).Value = 'Magnum'
$ws.Cells.Item(10, 10).Value = 100
$ws.Cells.Item(10, 11).Value = 40000
$ws.Cells.Item(10, 12).Value = 42000
$ws.Cells.Item(10, 13).Value = 41000
$ws.Cells.Item(10, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(10, 15).Value = 'Perú'
$ws.Cells.Item(10, 16).Value = 1640

# Row 11
$ws.Cells.Item(11, 4).Value = 44399
$ws.Cells.Item(11, 8).Value = 'Magnum'
$ws.Cells.Item(11, 10).Value = 100
$ws.Cells.Item(11, 11).Value = 20000
$ws.Cells.Item(11, 12).Value = 22000
$ws.Cells.Item(11, 13).Value = 21000
$ws.Cells.Item(11, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(11, 15).Value = 'Perú'
$ws.Cells.Item(11, 16).Value = 840

# Row 12
$ws.Cells.Item(12, 4).Value = 44435
$ws.Cells.Item(12, 8).Value = 'Magnum'
$ws.Cells.Item(12, 10).Value = 100
$ws.Cells.Item(12, 11).Value = 25000
$ws.Cells.Item(12, 12).Value = 26000
$ws.Cells.Item(12, 13).Value = 25500
$ws.Cells.Item(12, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(12, 15).Value = 'Perú'
$ws.Cells.Item(12, 16).Value = 1020

# Row 13
$ws.Cells.Item(13, 4).Value = 44203
$ws.Cells.Item(13, 8).Value = 'Magnum'
$ws.Cells.Item(13, 10).Value = 100
$ws.Cells.Item(13, 11).Value = 20000
$ws.Cells.Item(13, 12).Value = 22000
$ws.Cells.Item(13, 13).Value = 21000
$ws.Cells.Item(13, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(13, 15).Value = 'Región del Maule'
$ws.Cells.Item(13, 16).Value = 840

# Row 14
$ws.Cells.Item(14, 4).Value = 44475
$ws.Cells.Item(14, 8).Value = 'Magnum'
$ws.Cells.Item(14, 10).Value = 100
$ws.Cells.Item(14, 11).Value = 44000
$ws.Cells.Item(14, 12).Value = 45000
$ws.Cells.Item(14, 13).Value = 44500
$ws.Cells.Item(14, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(14, 15).Value = 'Perú'
$ws.Cells.Item(14, 16).Value = 1780

# Row 15
$ws.Cells.Item(15, 4).Value = 44244
$ws.Cells.Item(15, 8).Value = 'Magnum'
$ws.Cells.Item(15, 10).Value = 100
$ws.Cells.Item(15, 11).Value = 16000
$ws.Cells.Item(15, 12).Value = 18000
$ws.Cells.Item(15, 13).Value = 17000
$ws.Cells.Item(15, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(15, 15).Value = 'Región del Maule'
$ws.Cells.Item(15, 16).Value = 680

# Row 16
$ws.Cells.Item(16, 4).Value = 44167
$ws.Cells.Item(16, 8).Value = 'Sin especificar'
$ws.Cells.Item(16, 10).Value = 100
$ws.Cells.Item(16, 11).Value = 18000
$ws.Cells.Item(16, 12).Value = 19000
$ws.Cells.Item(16, 13).Value = 18500
$ws.Cells.Item(16, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(16, 15).Value = 'Región del Maule'
$ws.Cells.Item(16, 16).Value = 740

# Row 17
$ws.Cells.Item(17, 4).Value = 44321
$ws.Cells.Item(17, 8).Value = 'Magnum'
$ws.Cells.Item(17, 10).Value = 100
$ws.Cells.Item(17, 11).Value = 24000
$ws.Cells.Item(17, 12).Value = 25000
$ws.Cells.Item(17, 13).Value = 24500
$ws.Cells.Item(17, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(17, 15).Value = 'Región del Maule'
$ws.Cells.Item(17, 16).Value = 980

# Row 18
$ws.Cells.Item(18, 4).Value = 44265
$ws.Cells.Item(18, 8).Value = 'Magnum'
$ws.Cells.Item(18, 10).Value = 100
$ws.Cells.Item(18, 11).Value = 20000
$ws.Cells.Item(18, 12).Value = 22000
$ws.Cells.Item(18, 13).Value = 21000
$ws.Cells.Item(18, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(18, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(18, 16).Value = 840

# Row 19
$ws.Cells.Item(19, 4).Value = 44237
$ws.Cells.Item(19, 8).Value = 'Sin especificar'
$ws.Cells.Item(19, 10).Value = 100
$ws.Cells.Item(19, 11).Value = 20000
$ws.Cells.Item(19, 12).Value = 22000
$ws.Cells.Item(19, 13).Value = 21000
$ws.Cells.Item(19, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(19, 15).Value = 'Región del Maule'
$ws.Cells.Item(19, 16).Value = 840

# Row 20
$ws.Cells.Item(20, 4).Value = 44323
$ws.Cells.Item(20, 8).Value = 'Magnum'
$ws.Cells.Item(20, 10).Value = 100
$ws.Cells.Item(20, 11).Value = 20000
$ws.Cells.Item(20, 12).Value = 22000
$ws.Cells.Item(20, 13).Value = 21000
$ws.Cells.Item(20, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(20, 15).Value = 'Perú'
$ws.Cells.Item(20, 16).Value = 840

# Row 21
$ws.Cells.Item(21, 4).Value = 44342
$ws.Cells.Item(21, 8).Value = 'Magnum'
$ws.Cells.Item(21, 10).Value = 100
$ws.Cells.Item(21, 11).Value = 28000
$ws.Cells.Item(21, 12).Value = 30000
$ws.Cells.Item(21, 13).Value = 29000
$ws.Cells.Item(21, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(21, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(21, 16).Value = 1160

# Row 22
$ws.Cells.Item(22, 4).Value = 44363
$ws.Cells.Item(22, 8).Value = 'Magnum'
$ws.Cells.Item(22, 10).Value = 100
$ws.Cells.Item(22, 11).Value = 25000
$ws.Cells.Item(22, 12).Value = 26000
$ws.Cells.Item(22, 13).Value = 25500
$ws.Cells.Item(22, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(22, 15).Value = 'Perú'
$ws.Cells.Item(22, 16).Value = 1020

# Row 23
$ws.Cells.Item(23, 4).Value = 44441
$ws.Cells.Item(23, 8).Value = 'Magnum'
$ws.Cells.Item(23, 10).Value = 100
$ws.Cells.Item(23, 11).Value = 28000
$ws.Cells.Item(23, 12).Value = 29000
$ws.Cells.Item(23, 13).Value = 28500
$ws.Cells.Item(23, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(23, 15).Value = 'Perú'
$ws.Cells.Item(23, 16).Value = 1140

# Row 25
$ws.Cells.Item(25, 4).Value = 44384
$ws.Cells.Item(25, 8).Value = 'Sin especificar'
$ws.Cells.Item(25, 10).Value = 100
$ws.Cells.Item(25, 11).Value = 25000
$ws.Cells.Item(25, 12).Value = 26000
$ws.Cells.Item(25, 13).Value = 25500
$ws.Cells.Item(25, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(25, 15).Value = 'Perú'
$ws.Cells.Item(25, 16).Value = 1020

# Row 26
$ws.Cells.Item(26, 4).Value = 44272
$ws.Cells.Item(26, 8).Value = 'Magnum'
$ws.Cells.Item(26, 10).Value = 100
$ws.Cells.Item(26, 11).Value = 22000
$ws.Cells.Item(26, 12).Value = 24000
$ws.Cells.Item(26, 13).Value = 23000
$ws.Cells.Item(26, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(26, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(26, 16).Value = 920

# Row 27
$ws.Cells.Item(27, 4).Value = 44447
$ws.Cells.Item(27, 8).Value = 'Magnum'
$ws.Cells.Item(27, 10).Value = 100
$ws.Cells.Item(27, 11).Value = 37000
$ws.Cells.Item(27, 12).Value = 38000
$ws.Cells.Item(27, 13).Value = 37500
$ws.Cells.Item(27, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(27, 15).Value = 'Perú'
$ws.Cells.Item(27, 16).Value = 1500

# Row 28
$ws.Cells.Item(28, 4).Value = 44510
$ws.Cells.Item(28, 8).Value = 'Magnum'
$ws.Cells.Item(28, 10).Value = 100
$ws.Cells.Item(28, 11).Value = 35000
$ws.Cells.Item(28, 12).Value = 36000
$ws.Cells.Item(28, 13).Value = 35500
$ws.Cells.Item(28, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(28, 15).Value = 'Perú'
$ws.Cells.Item(28, 16).Value = 1420

# Row 29
$ws.Cells.Item(29, 4).Value = 44468
$ws.Cells.Item(29, 8).Value = 'Sin especificar'
$ws.Cells.Item(29, 10).Value = 100
$ws.Cells.Item(29, 11).Value = 31000
$ws.Cells.Item(29, 12).Value = 32000
$ws.Cells.Item(29, 13).Value = 31500
$ws.Cells.Item(29, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(29, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(29, 16).Value = 1260
